# Update cryptos list: refreshed prices and 1h volume percentages,
# and re-ranked WrappedEther/Polkadot (rows 13 and 14 swapped).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'27.102.38"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.74%  "
$ws.Range('E2').Style = 'Normal'
$ws.Range('D3').Value = "'1.849.24"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +0.26%  "
$ws.Range('E3').Style = 'Normal'
$ws.Range('D4').Value = "'1.018"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  +0.97%  "
$ws.Range('E4').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.85%  "
$ws.Range('E5').Style = 'Normal'
$ws.Range('D6').Value = "'309.89"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  +0.22%  "
$ws.Range('E6').Style = 'Normal'
$ws.Range('D7').Value = "'0.4774"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = "'  +2.01%  "
$ws.Range('E7').Style = 'Normal'
$ws.Range('D8').Value = "'0.3689"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = "'  +0.21%  "
$ws.Range('E8').Style = 'Normal'
$ws.Range('D9').Value = "'0.07254"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  +1.53%  "
$ws.Range('E9').Style = 'Normal'
$ws.Range('D10').Value = "'0.9308"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  +0.08%  "
$ws.Range('E10').Style = 'Normal'
$ws.Range('D11').Value = "'19.89"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +1.57%  "
$ws.Range('E11').Style = 'Normal'
$ws.Range('D12').Value = "'0.07790"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = "'  +1.19%  "
$ws.Range('E12').Style = 'Normal'
$ws.Range('B13').Value = "'WrappedEther"
$ws.Range('B13').Style = 'Normal'
$ws.Range('C13').Value = "'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range('C13').Style = 'Normal'
$ws.Range('D13').Value = "'1.857.93"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  -1.11%  "
$ws.Range('E13').Style = 'Normal'
$ws.Range('B14').Value = "'Polkadot"
$ws.Range('B14').Style = 'Normal'
$ws.Range('C14').Value = "'https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range('C14').Style = 'Normal'
$ws.Range('D14').Value = "'5.387"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  +1.90%  "
$ws.Range('E14').Style = 'Normal'
$ws.Range('D15').Value = "'6.478"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +1.25%  "
$ws.Range('E15').Style = 'Normal'
$ws.Range('E16').Value = "'  +1.46%  "
$ws.Range('E16').Style = 'Normal'
$ws.Range('D17').Value = "'1.017"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  +0.82%  "
$ws.Range('E17').Style = 'Normal'
$ws.Range('D18').Value = "'0.000008682"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.68%  "
$ws.Range('E18').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.77%  "
$ws.Range('E19').Style = 'Normal'
$ws.Range('D20').Value = "'27.110.81"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +0.65%  "
$ws.Range('E20').Style = 'Normal'
$ws.Range('D21').Value = "'14.58"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +1.41%  "
$ws.Range('E21').Style = 'Normal'
$ws.Range('D22').Value = "'5.065"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = "'  +0.74%  "
$ws.Range('E22').Style = 'Normal'
$ws.Range('E23').Value = "'  +0.04%  "
$ws.Range('E23').Style = 'Normal'
$ws.Range('D24').Value = "'1.937"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  +0.34%  "
$ws.Range('E24').Style = 'Normal'
$ws.Range('D25').Value = "'153.12"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  +0.48%  "
$ws.Range('E25').Style = 'Normal'
$ws.Range('E26').Value = "'  +0.57%  "
$ws.Range('E26').Style = 'Normal'
$ws.Range('D27').Value = "'1.986"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  -1.90%  "
$ws.Range('E27').Style = 'Normal'
$ws.Range('D28').Value = "'114.79"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = "'  +0.60%  "
$ws.Range('E28').Style = 'Normal'
$ws.Range('D29').Value = "'4.924"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  +0.77%  "
$ws.Range('E29').Style = 'Normal'
$ws.Range('D30').Value = "'0.08878"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.22%  "
$ws.Range('E30').Style = 'Normal'
$ws.Range('D31').Value = "'3.300"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +2.95%  "
$ws.Range('E31').Style = 'Normal'
$ws.Range('D32').Value = "'1.179"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  +0.46%  "
$ws.Range('E32').Style = 'Normal'
$ws.Range('D33').Value = "'4.524"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  +1.33%  "
$ws.Range('E33').Style = 'Normal'
$ws.Range('D34').Value = "'0.7365"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  -1.36%  "
$ws.Range('E34').Style = 'Normal'
$ws.Range('D35').Value = "'2.685"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  -4.09%  "
$ws.Range('E35').Style = 'Normal'
$ws.Range('D36').Value = "'1.117"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = "'  +3.09%  "
$ws.Range('E36').Style = 'Normal'
$ws.Range('E37').Value = "'  +2.09%  "
$ws.Range('E37').Style = 'Normal'
$ws.Range('D38').Value = "'0.05264"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  +1.44%  "
$ws.Range('E38').Style = 'Normal'
$ws.Range('E39').Value = "'  +0.42%  "
$ws.Range('E39').Style = 'Normal'
$ws.Range('D40').Value = "'0.5279"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = "'  +1.49%  "
$ws.Range('E40').Style = 'Normal'
$ws.Range('D41').Value = "'7.030"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  +1.93%  "
$ws.Range('E41').Style = 'Normal'
$ws.Range('D42').Value = "'0.1522"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = "'  +0.23%  "
$ws.Range('E42').Style = 'Normal'
$ws.Range('D43').Value = "'8.294"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  +1.91%  "
$ws.Range('E43').Style = 'Normal'
$ws.Range('D44').Value = "'10.61"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  +1.15%  "
$ws.Range('E44').Style = 'Normal'
$ws.Range('D45').Value = "'0.4736"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  +0.86%  "
$ws.Range('E45').Style = 'Normal'
$ws.Range('E46').Value = "'  +0.88%  "
$ws.Range('E46').Style = 'Normal'
$ws.Range('D47').Value = "'101.81"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.17%  "
$ws.Range('E47').Style = 'Normal'
$ws.Range('D48').Value = "'1.620"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.75%  "
$ws.Range('E48').Style = 'Normal'
$ws.Range('D49').Value = "'65.75"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  +1.43%  "
$ws.Range('E49').Style = 'Normal'
$ws.Range('D50').Value = "'0.06055"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  +0.44%  "
$ws.Range('E50').Style = 'Normal'
$ws.Range('D51').Value = "'0.8925"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  +0.34%  "
$ws.Range('E51').Style = 'Normal'
